$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "content" (old column E) for "n_message",
# pushing the old "content" column from E to F.
$ws.Columns(5).Insert()

# Header row
$ws.Range("E1").Value = "n_message"

# F1 needs the same header formatting (bold/centered/bordered) as the
# rest of row 1 - copy formats from a neighbouring header cell that
# already has it, then set the text.
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "content"

# Row 2: topic "Политический форум" (only columns A and B populated)
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Политический форум"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""

# Row 3: topic "Выборы в студсовет УлГТУ" (only columns A and B populated)
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Выборы в студсовет УлГТУ"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""

# Row 4: "ОТЗЫВ НАСТАВНИКА" feedback row
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "ОТЗЫВ НАСТАВНИКА"
$ws.Range("C4").Value = 522321184
$ws.Range("D4").Value = "Чебиняева Ирина Леонидовна"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "Чат-бот соответствует техническому заданию, все функции работают. Задание в рамках Обучения служением принято"

# The sheet now only needs 4 rows of data (plus header); drop the old row 5.
$ws.Rows(5).Delete()
